$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$null = $ws.Activate()

# Row 2 - ArriveNowCredentials
$ws.Range("B2").Value = "12_ArriveNowCredentials"
$ws.Range("A2").Value = "ArriveNowCredentials"
$ws.Range("C2").Value = "Shared"

# Column B for rows 3-9 (Asset names)
$ws.Range("B3").Value = "12_ArriveNowURL"
$ws.Range("B4").Value = "12_ArrivePortalURL"
$ws.Range("B5").Value = "12_ArriveTruckEntryURL"
$ws.Range("B6").Value = "12_EmailAddress"
$ws.Range("B7").Value = "12_GDriveReportFolder"
$ws.Range("B8").Value = "12_ToEmail"
$ws.Range("B9").Value = "12_CCEmail"

# Column A for rows 3-9 (Names)
$ws.Range("A3").Value = "ArriveNowURL"
$ws.Range("A4").Value = "ArrivePortalURL"
$ws.Range("A5").Value = "ArriveTruckEntryURL"
$ws.Range("A6").Value = "EmailAddress"
$ws.Range("A7").Value = "GDriveReportFolder"
$ws.Range("A8").Value = "ToEmail"
$ws.Range("A9").Value = "CCEmail"

# Column C for rows 3-9 (Folder)
$ws.Range("C3").Value = "Shared"
$ws.Range("C4").Value = "Shared"
$ws.Range("C5").Value = "Shared"
$ws.Range("C6").Value = "Shared"
$ws.Range("C7").Value = "Shared"
$ws.Range("C8").Value = "Shared"
$ws.Range("C9").Value = "Shared"

# Row 10 - ArriveCarrierSearchURL
$ws.Range("B10").Value = "12_ArriveCarrierSearchURL"
$ws.Range("A10").Value = "ArriveCarrierSearchURL"
$ws.Range("C10").Value = "Shared"

# Column D (Descriptions) in authoring order
$ws.Range("D2").Value = "Credentials to login into ArriveNow Portal"
$ws.Range("D3").Value = "URL for ArriveNow Portal"
$ws.Range("D4").Value = "URL for Arrive Portal"
$ws.Range("D5").Value = "URL for ArriveNow Truck Entry Portal"
$ws.Range("D10").Value = "URL for ArriveNow Carrier Search Portal"
$ws.Range("D6").Value = "Email account used to send and create report files"
$ws.Range("D7").Value = "ID for G Drive folder where reports are stored"
$ws.Range("D8").Value = "Email addresses where the emails are going to be sent to"
$ws.Range("D9").Value = "Email addresses copied to the emails that are going to be sent"

# Row 11 - ReportFileID
$ws.Range("B11").Value = "12_ReportFileID"
$ws.Range("A11").Value = "ReportFileID"
$ws.Range("C11").Value = "Shared"
$ws.Range("D11").Value = "ID for G Sheet used to report execution outputs"

# Final selection as left by the author
$null = $ws.Range("C18").Select()
